$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update sheet name (title reflects new "through" date)
$ws.Name = "Through 2021-12-17"

# Update the December label in column A, row 13
$ws.Range("A13").Value = "December (through 12-17)"

# Update December row (row 13) values
$ws.Range("B13").Value = 21
$ws.Range("C13").Value = 54
$ws.Range("D13").Value = 67
$ws.Range("E13").Value = 39
$ws.Range("F13").Value = 27
$ws.Range("G13").Value = 80
$ws.Range("H13").Value = 127

# Update Total row (row 14) values
$ws.Range("B14").Value = 312
$ws.Range("C14").Value = 617
$ws.Range("D14").Value = 888
$ws.Range("E14").Value = 721
$ws.Range("F14").Value = 561
$ws.Range("G14").Value = 1344
$ws.Range("H14").Value = 1770
